$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 22:20"

# Update country names (column A) to reflect new sort/reorder of countries list
$ws.Range("A91").Value = "San Marino"
$ws.Range("A92").Value = "Camerun"
$ws.Range("A93").Value = "Vietnam"
$ws.Range("A94").Value = "Cuba"
$ws.Range("A95").Value = "Oman"
$ws.Range("A96").Value = "Afganistan"
$ws.Range("A166").Value = "Benin"
$ws.Range("A167").Value = "Dominica"
$ws.Range("A168").Value = "Guyana"
$ws.Range("A169").Value = "Curazao"
$ws.Range("A171").Value = "Seychelles"
$ws.Range("A172").Value = "Laos"
$ws.Range("A174").Value = "Mozambique"
$ws.Range("A175").Value = "Siria"
$ws.Range("A176").Value = "Groenlandia"
$ws.Range("A178").Value = "Granada"
$ws.Range("A182").Value = "Republica del Chad"
$ws.Range("A183").Value = "Antigua y Barbuda"
$ws.Range("A186").Value = "Liberia"
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("A188").Value = "Santa Sede"
$ws.Range("A192").Value = "Montserrat"
$ws.Range("A193").Value = "Fiyi"
$ws.Range("A194").Value = "Nepal"
$ws.Range("A196").Value = "Somalia"
$ws.Range("A198").Value = "Botsuana"
$ws.Range("A199").Value = "Gambia"
$ws.Range("A200").Value = "Belice"
$ws.Range("A201").Value = "Islas Virgenes Britanicas"
$ws.Range("A202").Value = "Republica de Africa Central"
$ws.Range("A203").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A204").Value = "Anguila"
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Papua Nueva Guinea"

# Update statistic values (columns B-H) with refreshed data
$ws.Range("B4").Value = 210714
$ws.Range("C4").Value = 22184
$ws.Range("E4").Value = 197212
$ws.Range("G4").Value = 644
$ws.Range("H4").Value = 4697
$ws.Range("B8").Value = 77779
$ws.Range("C8").Value = 5971
$ws.Range("E8").Value = 58170
$ws.Range("G8").Value = 134
$ws.Range("H8").Value = 909
$ws.Range("B38").Value = 1998
$ws.Range("C38").Value = 601
$ws.Range("E38").Value = 1792
$ws.Range("G38").Value = 23
$ws.Range("H38").Value = 58
$ws.Range("B72").Value = 459
$ws.Range("C72").Value = 39
$ws.Range("E72").Value = 427
$ws.Range("B91").Value = 236
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 13
$ws.Range("E91").Value = 197
$ws.Range("F91").Value = 16
$ws.Range("H91").Value = 26
$ws.Range("B92").Value = 233
$ws.Range("C92").Value = 40
$ws.Range("D92").Value = 10
$ws.Range("E92").Value = 217
$ws.Range("F92").Value = 0
$ws.Range("H92").Value = 6
$ws.Range("B93").Value = 218
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 63
$ws.Range("E93").Value = 155
$ws.Range("F93").Value = 3
$ws.Range("H93").Value = 0
$ws.Range("B94").Value = 212
$ws.Range("C94").Value = 26
$ws.Range("D94").Value = 12
$ws.Range("E94").Value = 194
$ws.Range("H94").Value = 6
$ws.Range("B95").Value = 210
$ws.Range("C95").Value = 18
$ws.Range("D95").Value = 34
$ws.Range("E95").Value = 175
$ws.Range("H95").Value = 1
$ws.Range("B96").Value = 196
$ws.Range("C96").Value = 22
$ws.Range("D96").Value = 5
$ws.Range("E96").Value = 187
$ws.Range("F96").Value = 0
$ws.Range("H96").Value = 4
$ws.Range("B120").Value = 90
$ws.Range("C120").Value = 3
$ws.Range("E120").Value = 84
$ws.Range("B166").Value = 13
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 1
$ws.Range("E167").Value = 12
$ws.Range("H167").Value = 0
$ws.Range("B168").Value = 12
$ws.Range("D168").Value = 0
$ws.Range("E168").Value = 10
$ws.Range("H168").Value = 2
$ws.Range("B169").Value = 11
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 3
$ws.Range("E169").Value = 7
$ws.Range("H169").Value = 1
$ws.Range("C171").Value = 0
$ws.Range("C172").Value = 1
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 10
$ws.Range("C174").Value = 2
$ws.Range("E174").Value = 10
$ws.Range("H174").Value = 0
$ws.Range("D175").Value = 0
$ws.Range("H175").Value = 2
$ws.Range("B176").Value = 10
$ws.Range("D176").Value = 2
$ws.Range("E176").Value = 8
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 9
$ws.Range("C186").Value = 3
$ws.Range("C187").Value = 1
$ws.Range("C188").Value = 0
$ws.Range("C203").Value = 2
$ws.Range("C204").Value = 0
